$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain-text strings in the source feed (e.g.
# "218.68" or "26.270.95"), not numbers. Assigning such a string straight to
# .Value lets Excel auto-convert simple decimals (e.g. "218.68") into a
# numeric cell, which would not match the expected text cells. For every
# Price cell being updated, force the text number format first so the value
# is kept as text, then restore the default "Normal" style afterwards so no
# visible formatting change is introduced.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.270.95"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.680.48"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "218.68"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "0.5275"
$ws.Range("E6").Value = "  +2.69%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.2702"
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").Value = "0.06433"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").Value = "22.08"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("D11").Value = "0.07507"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.687.68"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.552"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "0.5815"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "0.000008479"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").Value = "64.35"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "26.319.64"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "4.926"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").Value = "189.56"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "6.209"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "1.008"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "144.82"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "7.735"
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("D26").Value = "0.1237"
$ws.Range("E26").Value = "  +4.44%  "
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").Value = "0.06648"
$ws.Range("E28").Value = "  +11.75%  "
$ws.Range("D29").Value = "1.359"
$ws.Range("E29").Value = "  +5.96%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "1.663"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").Value = "1.027"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("E35").Value = "  +3.02%  "
$ws.Range("D37").Value = "2.705"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").Value = "6.397"
$ws.Range("E38").Value = "  +5.11%  "
$ws.Range("D39").Value = "1.108.30"
$ws.Range("E39").Value = "  +2.64%  "
$ws.Range("D40").Value = "0.01624"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").Value = "0.8769"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").Value = "1.014"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").Value = "100.48"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "1.828.92"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("D45").Value = "0.00000000112"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").Value = "56.92"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").Value = "8.188"
$ws.Range("E47").Value = "  +1.67%  "
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("D50").Value = "0.4303"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "6.053"
$ws.Range("E51").Value = "  +2.84%  "

$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
